# Generate Report for Handoff
# Update the "Latest Handoff Datetime" (column D, row 5 -- the
# ebf404d6-2044-42ed-87dd-412926134f40.md entry) on each locale sheet with
# the timestamp of the newly generated handoff.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-03-10 09:04:08"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-03-10 09:04:17"
